# Adds a new "Player Info" sheet (as the first/leftmost sheet) with basic
# player identity fields, and reshapes the "ODI Batting" sheet's
# MATCH_CARD_LINK column into a plain MATCH_CODE number column.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert the new "Player Info" sheet. Worksheets.Add() drops the new
#    sheet in front of the active sheet, i.e. at the very start of the
#    workbook - which is exactly where it needs to land.
# ------------------------------------------------------------------
$info = $wb.Worksheets.Add()
$info.Name = "Player Info"

# ------------------------------------------------------------------
# 2. Grab a handle to the original "ODI Batting" sheet. This MUST be
#    fetched after the insert/rename above - earlier sheet references
#    can resolve by position and go stale once the sheet collection is
#    reshuffled.
# ------------------------------------------------------------------
$batting = $wb.Worksheets.Item("ODI Batting")

function Set-TextValue {
    param($range, [string]$text)
    # Force the cell to be written as TEXT even when the content looks
    # numeric (e.g. "7129"), matching the source data pipeline, which
    # always emits inlineStr cells. Reset back to the default "Normal"
    # style afterwards so no stray number-format styling is left on the
    # cell (only its text-ness is kept).
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

function Set-HeaderCell {
    param($range, [string]$text)
    $range.Value = $text
    $range.Font.Bold = $true
    $range.HorizontalAlignment = -4108   # xlCenter
    $range.VerticalAlignment = -4160     # xlTop
    $range.Borders.LineStyle = 1         # xlContinuous (thin box border)
}

# ------------------------------------------------------------------
# 3. Header row for Player Info.
# ------------------------------------------------------------------
Set-HeaderCell $info.Range("A1") "ID"
Set-HeaderCell $info.Range("B1") "NAME"
Set-HeaderCell $info.Range("C1") "BATTING_HAND"
Set-HeaderCell $info.Range("D1") "BOWL_STYLE"

# ------------------------------------------------------------------
# 4. Data row for Player Info.
# ------------------------------------------------------------------
Set-TextValue $info.Range("A2") "7129"
Set-TextValue $info.Range("B2") "Tomas Scott Sabater Mackintosh"
Set-TextValue $info.Range("C2") "Right Handed"
Set-TextValue $info.Range("D2") "Does Not Bowl | Unknown"

$info.Range("A1").Select()

# ------------------------------------------------------------------
# 5. Rework the ODI Batting sheet: MATCH_CARD_LINK -> MATCH_CODE, and
#    the scorecard URLs collapse down to just the trailing MatchCode
#    number.
# ------------------------------------------------------------------
$batting.Range("D1").Value = "MATCH_CODE"

Set-TextValue $batting.Range("D2") "4703"
Set-TextValue $batting.Range("D3") "4705"
Set-TextValue $batting.Range("D4") "4706"

Write-Output "Player Info sheet added; ODI Batting MATCH_CODE column updated."
